# Update crypto price/volume figures for the Mon Sep 18 13:29:59 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.386.38"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "1.666.04"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.254"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "1.899.64"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "1.669.71"
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.26%  "
$ws.Range("D17").Value = "27.364.73"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "224.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.78%  "
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.95%  "
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("E27").Value = "  +3.66%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  +2.73%  "
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  +2.65%  "
$ws.Range("D35").Value = "1.270.18"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.832"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.813"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").Value = "1.811.26"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("E44").Value = "  -4.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0983"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  -0.01%  "
